# Applies the cryptos.xlsx price-refresh diff (Sat Aug 12 23:32:11 UTC 2023 GitHub Actions run).
# For cells whose new text is a bare decimal number (e.g. "240.57"), a leading
# apostrophe forces Excel to keep storing it as literal text (matching the source
# inline-string cells) instead of silently converting it to a Number value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "29.434.01"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "1.849.87"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 (BNB)
$ws.Range("D5").Value = "'240.57"
$ws.Range("E5").Value = "  +0.43%  "

# Row 6 (XRP)
$ws.Range("D6").Value = "'0.6276"
$ws.Range("E6").Value = "  -0.75%  "

# Row 7 (USDC)
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 (Dogecoin)
$ws.Range("D8").Value = "'0.07682"
$ws.Range("E8").Value = "  +1.50%  "

# Row 9 (Cardano)
$ws.Range("D9").Value = "'0.2921"
$ws.Range("E9").Value = "  -0.46%  "

# Row 10 (Solana)
$ws.Range("E10").Value = "  +1.34%  "

# Row 12 (WrappedEther)
$ws.Range("D12").Value = "1.853.64"
$ws.Range("E12").Value = "  -0.28%  "

# Row 13 (Polkadot)
$ws.Range("D13").Value = "'5.039"
$ws.Range("E13").Value = "  +0.69%  "

# Row 14 (ShibaInu)
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6818"
$ws.Range("E14").Value = "  +0.19%  "

# Row 15 (Polygon)
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001077"
$ws.Range("E15").Value = "  +3.35%  "

# Row 16 (Litecoin)
$ws.Range("D16").Value = "'83.53"
$ws.Range("E16").Value = "  +0.13%  "

# Row 17 (WrappedliquidstakedEther2.0)
$ws.Range("D17").Value = "2.114.84"
$ws.Range("E17").Value = "  +0.05%  "

# Row 18 (Uniswap)
$ws.Range("D18").Value = "'6.219"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19 (WrappedBTC)
$ws.Range("D19").Value = "29.462.70"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20 (BitcoinCash)
$ws.Range("D20").Value = "'228.85"
$ws.Range("E20").Value = "  -0.15%  "

# Row 21 (Avalanche)
$ws.Range("D21").Value = "'12.40"
$ws.Range("E21").Value = "  -0.30%  "

# Row 22 (Dai)
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23 (Chainlink)
$ws.Range("E23").Value = "  -0.46%  "

# Row 24 (BinanceUSD)
$ws.Range("E24").Value = "  +0.04%  "

# Row 25 (Monero)
$ws.Range("D25").Value = "'157.68"
$ws.Range("E25").Value = "  +0.58%  "

# Row 26 (Stellar)
$ws.Range("E26").Value = "  -1.08%  "

# Row 27 (Cosmos)
$ws.Range("D27").Value = "'8.417"
$ws.Range("E27").Value = "  +0.73%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").Value = "'17.74"
$ws.Range("E28").Value = "  +0.74%  "

# Row 29 (Toncoin)
$ws.Range("D29").Value = "'1.367"
$ws.Range("E29").Value = "  +5.00%  "

# Row 30 (PancakeSwap)
$ws.Range("D30").Value = "'1.463"
$ws.Range("E30").Value = "  +0.29%  "

# Row 31 (Hedera)
$ws.Range("D31").Value = "'0.05633"
$ws.Range("E31").Value = "  -0.30%  "

# Row 32 (Filecoin)
$ws.Range("D32").Value = "'4.128"
$ws.Range("E32").Value = "  +0.66%  "

# Row 33 (InternetComputer(DFINITY))
$ws.Range("D33").Value = "'4.060"
$ws.Range("E33").Value = "  +0.71%  "

# Row 34 (LidoDAOToken)
$ws.Range("D34").Value = "'1.846"
$ws.Range("E34").Value = "  -0.13%  "

# Row 35 (ARBITRUM)
$ws.Range("E35").Value = "  +0.45%  "

# Row 36 (ImmutableX)
$ws.Range("D36").Value = "'0.7070"
$ws.Range("E36").Value = "  -0.49%  "

# Row 37 (HuobiToken)
$ws.Range("D37").Value = "'2.597"
$ws.Range("E37").Value = "  +0.16%  "

# Row 38 (Maker)
$ws.Range("D38").Value = "1.226.17"
$ws.Range("E38").Value = "  -1.97%  "

# Row 39 (VeChain)
$ws.Range("D39").Value = "'0.01794"

# Row 40 (MXToken)
$ws.Range("D40").Value = "'2.758"
$ws.Range("E40").Value = "  -0.48%  "

# Row 41 (FraxShare)
$ws.Range("D41").Value = "'6.460"
$ws.Range("E41").Value = "  +1.25%  "

# Row 42 (TrustWalletToken)
$ws.Range("D42").Value = "'0.9036"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43 (PaxDollar)
$ws.Range("E43").Value = "  +0.07%  "

# Row 44 (RocketPoolETH)
$ws.Range("D44").Value = "2.021.32"
$ws.Range("E44").Value = "  -0.06%  "

# Row 45 (Quant)
$ws.Range("D45").Value = "'102.01"
$ws.Range("E45").Value = "  +0.21%  "

# Row 46 (Aave)
$ws.Range("D46").Value = "'66.16"
$ws.Range("E46").Value = "  +0.35%  "

# Row 47 (Aptos)
$ws.Range("D47").Value = "'7.187"
$ws.Range("E47").Value = "  +1.28%  "

# Row 48 (BabyDogeCoin)
$ws.Range("D48").Value = "'0.00000000117"
$ws.Range("E48").Value = "  -0.94%  "

# Row 49 (TheSandbox)
$ws.Range("D49").Value = "'0.4023"
$ws.Range("E49").Value = "  +0.46%  "

# Row 50 (Algorand)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.055"
$ws.Range("E50").Value = "  +0.97%  "

# Row 51 (EnergySwap)
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.1157"
$ws.Range("E51").Value = "  +2.97%  "
